# daily auto push: 2026-02-20 14:04 UTC
# A new reading was logged for 2026/02/20 (金) at hour 20. This pushes a
# new row into the existing "2026/02/20" block (row 838) and shifts every
# following row down by one, so the table grows from A1:D879 to A1:D880.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row right before the old row 838, pushing rows 838-879
# down to 839-880 (this also grows the sheet's used range automatically).
$ws.Rows("838:838").Insert()

# Column A holds a date-like literal string (e.g. "2026/02/20"), not a
# real Excel date. Format the cell as Text first so the slash-separated
# value isn't auto-coerced into a date serial number, then reset the
# style back to Normal so the new row matches the formatting (no
# explicit style) of every other data row in the table.
$ws.Cells.Item(838, 1).NumberFormat = "@"
$ws.Cells.Item(838, 1).Value = "2026/02/20"
$ws.Cells.Item(838, 1).Style = "Normal"

$ws.Cells.Item(838, 2).Value = "金"
$ws.Cells.Item(838, 3).Value = 20
$ws.Cells.Item(838, 4).Value = 201
